$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 91 (pushes existing rows 91.. down by one,
# carrying the row-91 cell formatting, e.g. the date number format in D).
$ws.Rows(91).Insert()
$ws.Range("A91").Value = 4
$ws.Range("B91").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C91").Value = "Los Lagos"
$ws.Range("D91").Value = 44473
$ws.Range("E91").Value = 10
$ws.Range("F91").Value = 100112003
$ws.Range("G91").Value = "Ajo"
$ws.Range("H91").Value = "Chino"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 80
$ws.Range("K91").Value = 20000
$ws.Range("L91").Value = 20000
$ws.Range("M91").Value = 20000
$ws.Range("N91").Value = "$/caja 10 kilos"
$ws.Range("O91").Value = "China"
$ws.Range("P91").Value = 2000
$ws.Range("Q91").Value = 10
$ws.Range("R91").Value = "Hortaliza"

# Insert a second new data row at row 135 (after the first insertion has
# already shifted everything down by one).
$ws.Rows(135).Insert()
$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44476
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = 100112003
$ws.Range("G135").Value = "Ajo"
$ws.Range("H135").Value = "Chino"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 100
$ws.Range("K135").Value = 20000
$ws.Range("L135").Value = 20000
$ws.Range("M135").Value = 20000
$ws.Range("N135").Value = "$/caja 10 kilos"
$ws.Range("O135").Value = "China"
$ws.Range("P135").Value = 2000
$ws.Range("Q135").Value = 10
$ws.Range("R135").Value = "Hortaliza"
